# Repull data, push all data, mean calculation
# Update column F (dSF) values on Sheet1 to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 2
    4  = 4
    5  = -1
    6  = 3
    7  = 3
    9  = 6
    10 = -6
    11 = -1
    12 = 1
    13 = 10
    14 = -1
    15 = -1
    16 = 3
    17 = -1
    18 = 2
    19 = 0
    20 = 3
    21 = 2
    22 = 3
    24 = -4
    25 = -4
    26 = 4
    27 = -3
    28 = 2
    29 = 1
    31 = 3
    32 = -5
    33 = -7
    34 = 1
    35 = 6
    36 = 1
    37 = 4
    38 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
